$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.644.09'
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").Value = '2.334.40'
$ws.Range("E3").Value = '  +1.76%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.61%  '

$ws.Range("E7").Value = '  +1.02%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +2.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0919'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.64%  '

$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("E13").Value = '  -0.90%  '

$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").Value = '2.684.54'
$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("D17").Value = '2.329.61'
$ws.Range("E17").Value = '  +1.15%  '

$ws.Range("D18").Value = '43.545.20'
$ws.Range("E18").Value = '  +1.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.31%  '

$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.30%  '

$ws.Range("E25").Value = '  +3.16%  '

$ws.Range("E26").Value = '  -0.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '

$ws.Range("E29").Value = '  -2.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.05%  '

$ws.Range("E31").Value = '  +1.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.50'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.54%  '

$ws.Range("E33").Value = '  +1.63%  '

$ws.Range("E34").Value = '  +9.85%  '

$ws.Range("E35").Value = '  +0.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.78'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.46%  '

$ws.Range("E37").Value = '  -1.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0366'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.71%  '

$ws.Range("E40").Value = '  +6.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.08%  '

$ws.Range("E43").Value = '  +3.03%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.236'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.77%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '114.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("D48").Value = '1.669.62'
$ws.Range("E48").Value = '  -2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '77.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.30%  '

$ws.Range("E51").Value = '  +1.97%  '
